$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the Instance value in row 2 from "bob" to "automation"
$ws.Range("D2").Value = "automation"

# Move the selection / active cell to D2 (single cell, matching sort/settings update)
$ws.Range("D2").Select()
